# Apply updated cryptocurrency price/volume data to sheet1 (Coin list)
# Generated from the canonical-OOXML diff describing the May 27 2023 GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings that LOOK like plain numbers (e.g. "1.012", "0.07210").
# Force those specific cells to Text format *before* assigning the value so Excel keeps the
# exact original string (incl. trailing zeros) instead of silently converting it to a number.
$numericLookingPriceCells = @(
    "D4",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Column D = Price, Column E = Volume(1h) change percentage
$ws.Range("D2").Value = "26.965.29"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.843.81"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "308.66"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "0.4767"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("D8").Value = "0.3674"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "0.07210"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "0.9301"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "19.78"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "0.07723"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.865.97"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "5.390"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "6.441"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "88.84"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "0.000008655"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "26.987.02"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "5.068"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "1.940"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").Value = "152.48"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "18.18"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").Value = "2.011"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "114.29"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").Value = "0.08855"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "3.295"
$ws.Range("E31").Value = "  +4.26%  "
$ws.Range("D32").Value = "1.175"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "0.7396"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "4.494"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "2.697"
$ws.Range("E35").Value = "  -5.39%  "
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").Value = "0.01959"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "0.05253"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("D39").Value = "2.965"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "0.5241"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").Value = "7.001"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "0.1510"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "8.294"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "10.53"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "0.4734"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "1.011"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "101.72"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "65.61"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").Value = "0.06070"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "0.8880"
$ws.Range("E51").Value = "  +3.19%  "
